$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: [ref, newValue, forceText]
$updates = @(
    @('D2', '43.006.26', $true),
    @('E2', '  -0.34%  ', $true),
    @('D3', '2.553.75', $true),
    @('E3', '  +0.11%  ', $true),
    @('E4', '  +0.05%  ', $true),
    @('D5', '304.57', $true),
    @('E5', '  +2.02%  ', $true),
    @('D6', '98.49', $true),
    @('E6', '  +4.32%  ', $true),
    @('E7', '  +0.01%  ', $true),
    @('E8', '  +0.04%  ', $true),
    @('D9', '0.549', $true),
    @('E9', '  -0.52%  ', $true),
    @('D10', '37.13', $true),
    @('E10', '  +2.87%  ', $true),
    @('D11', '0.0816', $true),
    @('E11', '  +0.55%  ', $true),
    @('E12', '  -0.03%  ', $true),
    @('E13', '  +6.79%  ', $true),
    @('B14', 'WrappedliquidstakedEther2.0', $false),
    @('C14', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', $false),
    @('D14', '2.948.15', $true),
    @('E14', '  +0.23%  ', $true),
    @('B15', 'WrappedEther', $false),
    @('C15', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', $false),
    @('D15', '2.637.46', $true),
    @('E15', '  +2.98%  ', $true),
    @('B16', 'Chainlink', $false),
    @('C16', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', $false),
    @('D16', '14.92', $true),
    @('E16', '  +4.98%  ', $true),
    @('B17', 'Polygon', $false),
    @('C17', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', $false),
    @('D17', '0.884', $true),
    @('E17', '  +1.13%  ', $true),
    @('B18', 'WrappedBTC', $false),
    @('C18', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', $false),
    @('D18', '43.068.97', $true),
    @('E18', '  -0.19%  ', $true),
    @('B19', 'InternetComputer(DFINITY)', $false),
    @('C19', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', $false),
    @('D19', '13.78', $true),
    @('E19', '  +5.85%  ', $true),
    @('B20', 'ShibaInu', $false),
    @('C20', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', $false),
    @('D20', '0.0₃0990', $true),
    @('E20', '  +0.62%  ', $true),
    @('B21', 'Uniswap', $false),
    @('C21', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', $false),
    @('D21', '6.63', $true),
    @('E21', '  -0.42%  ', $true),
    @('B22', 'Litecoin', $false),
    @('C22', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', $false),
    @('D22', '71.91', $true),
    @('E22', '  -0.35%  ', $true),
    @('B23', 'BitcoinCash', $false),
    @('C23', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', $false),
    @('D23', '256.29', $true),
    @('E23', '  -1.72%  ', $true),
    @('B24', 'PancakeSwap', $false),
    @('C24', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', $false),
    @('D24', '2.97', $true),
    @('E24', '  +1.97%  ', $true),
    @('B25', 'ImmutableX', $false),
    @('C25', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', $false),
    @('D25', '2.10', $true),
    @('E25', '  -2.33%  ', $true),
    @('B26', 'EthereumClassic', $false),
    @('C26', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', $false),
    @('D26', '28.09', $true),
    @('E26', '  -4.57%  ', $true),
    @('B27', 'Dai', $false),
    @('C27', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', $false),
    @('D27', '0.999', $true),
    @('E27', '  +0.00%  ', $true),
    @('B28', 'Cosmos', $false),
    @('C28', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', $false),
    @('D28', '10.17', $true),
    @('E28', '  +1.06%  ', $true),
    @('B29', 'InjectiveProtocol', $false),
    @('C29', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', $false),
    @('D29', '38.09', $true),
    @('E29', '  +2.89%  ', $true),
    @('B30', 'Toncoin', $false),
    @('C30', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', $false),
    @('D30', '2.10', $true),
    @('E30', '  -1.31%  ', $true),
    @('B31', 'Filecoin', $false),
    @('C31', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', $false),
    @('D31', '6.05', $true),
    @('E31', '  +0.57%  ', $true),
    @('B32', 'Monero', $false),
    @('C32', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', $false),
    @('D32', '158.72', $true),
    @('E32', '  +2.21%  ', $true),
    @('B33', 'ARBITRUM', $false),
    @('C33', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', $false),
    @('D33', '2.16', $true),
    @('E33', '  -0.28%  ', $true),
    @('B34', 'WEMIXToken', $false),
    @('C34', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', $false),
    @('D34', '2.76', $true),
    @('E34', '  +0.41%  ', $true),
    @('B35', 'Hedera', $false),
    @('C35', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', $false),
    @('D35', '0.0811', $true),
    @('E35', '  +1.27%  ', $true),
    @('B36', 'LidoDAOToken', $false),
    @('C36', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', $false),
    @('D36', '3.33', $true),
    @('E36', '  -1.85%  ', $true),
    @('B37', 'Celestia', $false),
    @('C37', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia', $false),
    @('D37', '18.98', $true),
    @('E37', '  +14.30%  ', $true),
    @('B38', 'EnergySwap', $false),
    @('C38', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', $false),
    @('D38', '25.77', $true),
    @('E38', '  +10.86%  ', $true),
    @('B39', 'Kaspa', $false),
    @('C39', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', $false),
    @('D39', '0.115', $true),
    @('E39', '  -0.74%  ', $true),
    @('B40', 'Stellar', $false),
    @('C40', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', $false),
    @('D40', '0.119', $true),
    @('E40', '  -0.49%  ', $true),
    @('B41', 'NEARProtocol', $false),
    @('C41', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', $false),
    @('D41', '3.46', $true),
    @('E41', '  -0.73%  ', $true),
    @('D42', '3.88', $true),
    @('E42', '  -0.44%  ', $true),
    @('B43', 'ApeXProtocol', $false),
    @('C43', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', $false),
    @('D43', '2.06', $true),
    @('E43', '  +29.54%  ', $true),
    @('B44', 'Maker', $false),
    @('C44', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', $false),
    @('D44', '2.099.86', $true),
    @('E44', '  +1.23%  ', $true),
    @('B45', 'VeChain', $false),
    @('C45', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', $false),
    @('D45', '0.0306', $true),
    @('E45', '  -2.44%  ', $true),
    @('B46', 'FirstDigitalUSD', $false),
    @('C46', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', $false),
    @('D46', '0.999', $true),
    @('E46', '  +0.07%  ', $true),
    @('B47', 'BitcoinSV', $false),
    @('C47', 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv', $false),
    @('D47', '87.24', $true),
    @('E47', '  +1.39%  ', $true),
    @('B48', 'FraxShare', $false),
    @('C48', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', $false),
    @('D48', '9.03', $true),
    @('E48', '  +3.14%  ', $true),
    @('B49', 'ordi', $false),
    @('C49', 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi', $false),
    @('D49', '75.55', $true),
    @('E49', '  +8.79%  ', $true),
    @('B50', 'RocketPoolETH', $false),
    @('C50', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', $false),
    @('D50', '2.804.59', $true),
    @('E50', '  +0.25%  ', $true),
    @('B51', 'Aave', $false),
    @('C51', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', $false),
    @('D51', '103.85', $true),
    @('E51', '  -0.55%  ', $true),
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    if ($forceText) {
        $ws.Range($ref).Formula = "'" + $val
        $ws.Range($ref).Style = "Normal"
    } else {
        $ws.Range($ref).Value2 = $val
    }
}